$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update timing measurements (B5:B7) ---
$ws.Range("B5").Value = 0.0007128715515136719
$ws.Range("B6").Value = 0.0005338191986083984
$ws.Range("B7").Value = 0.00877690315246582

# --- Convert tuple-style text "(x, y)" to list-style text "[x, y]" in the
#     embedding/permutation rows ---
$ws.Range("A8").Value = "[[2, 3], [1, 3], [3, 2], [3, 0], [2, 0], [2, 2], [3, 1], [0, 3], [0, 2], [0, 0], [2, 1], [1, 2], [1, 0], [1, 1], [0, 1]]"
$ws.Range("A50").Value = "[[2, 3], [3, 3], [3, 2], [3, 0], [2, 0], [2, 1], [1, 1], [1, 2], [0, 2], [0, 3], [1, 3], [0, 1], [2, 2], [1, 0], [3, 1]]"
$ws.Range("A111").Value = "[[2, 3], [1, 3], [3, 2], [3, 1], [2, 0], [1, 0], [0, 2], [0, 1], [0, 0], [1, 1], [2, 1], [2, 2], [3, 3], [1, 2], [0, 3]]"
$ws.Range("A177").Value = "[[3, 2], [3, 1], [2, 0], [1, 0], [0, 1], [0, 0], [1, 1], [2, 1], [3, 0], [2, 2], [3, 3], [1, 2], [0, 3], [0, 2], [1, 3]]"
$ws.Range("A236").Value = "[[1, 1], [0, 2], [0, 1], [1, 2], [0, 0], [1, 3], [1, 0], [3, 0], [2, 1], [2, 3], [2, 0], [3, 2], [0, 3], [2, 2], [3, 1]]"

# --- Insert a new summary row before the old row 284 ("Movement times"),
#     pushing the footer block (rows 284-288) down to 285-289 ---
$ws.Rows.Item(284).Insert()

# Populate the newly inserted row 284 with the move_fidelity result
$ws.Range("A284").Value = "move_fidelity"
$ws.Range("B284").Value = 0.997676961296731

# Update the "total time:" value, now on row 288
$ws.Range("B288").Value = 0.03292083740234375
